$d = $word.ActiveDocument

# XML fragment used with Range.InsertXML to create a brand-new run that has
# NO run properties (no <w:rPr>), matching the target diff which introduces
# a run consisting solely of 16 spaces with no formatting.
$spacesXmlTemplate = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">{0}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$sixteenSpaces = "                "
$spacesXml = $spacesXmlTemplate -f $sixteenSpaces

# Paragraphs whose leading <w:tab/> must simply be deleted (formatting of the
# remaining run/text is otherwise untouched).
$simpleTabParagraphs = @(
    "Rinse fully ripe figs.",
    "Cuts the figs in half.",
    "Lay them on a wire or wooden rack covered with cheesecloth.",
    "Cover the figs with cheesecloth.",
    "Place the rack in full sunlight during the day.",
    "Return the figs to the sun for 2 to 3 days.",
    "Store the dried figs in airtight containers in the fridge or freezer."
)

foreach ($snippet in $simpleTabParagraphs) {
    $search = $d.Content
    $found = $search.Find.Execute("^t" + $snippet, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "WARNING: could not find tab + '$snippet'"
        continue
    }
    # $search now covers "<tab>Text..."; shrink it to just the leading tab
    # character and delete that character only.
    $tabRange = $d.Range($search.Start, $search.Start + 1)
    $tabRange.Text = ""
}

# Special case: "Currants figs" paragraph. The original tab is replaced by a
# brand new, separate run containing 16 spaces (with no run formatting),
# while the original run (which keeps its en-US language formatting) loses
# only the <w:tab/> and keeps its text.
$search = $d.Content
$found = $search.Find.Execute("^tCurrants figs", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    Write-Host "WARNING: could not find tab + 'Currants figs'"
} else {
    $tabStart = $search.Start
    $tabRange = $d.Range($tabStart, $tabStart + 1)
    # Delete the tab character first.
    $tabRange.Text = ""
    # Insert the new, unformatted 16-space run at the exact same position,
    # using a collapsed range so InsertXML inserts in place rather than at
    # the end of the paragraph.
    $insertionPoint = $d.Range($tabStart, $tabStart)
    $insertionPoint.InsertXML($spacesXml)
}
